# Updated PDF Merger to latest version
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the score/total headers to intensity-based wording
$ws.Range("J1").Value = "Sun Aspect Intensity"
$ws.Range("K1").Value = "Moon Aspect Intensity"
$ws.Range("L1").Value = "Asc Aspect Intensity"
$ws.Range("M1").Value = "Aspect intensity total"

# Add a running "intensity total" formula column in M for data rows 2-91,
# each referencing the aspect-score row immediately above it.
for ($row = 2; $row -le 91; $row++) {
    $prevRow = $row - 1
    $ws.Range("M$row").Formula = "=SUM(J${prevRow}:L${prevRow})"
}
